$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.142637666666667
$ws.Range("H2").Value = 6.427913
$ws.Range("I2").Value = 0.05953067687027142
$ws.Range("J2").Value = 0.05953067687027141
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.236283
$ws.Range("N2").Value = 0.708849
$ws.Range("O2").Value = 0.001461516295904947
$ws.Range("P2").Value = 0.001461516295904947
$ws.Range("Q2").Value = 0.5062688557930001
$ws.Range("R2").Value = 4.556419702137
$ws.Range("S2").Value = 0.00008700505435215339
$ws.Range("T2").Value = 0.00008700505435215336
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.142637666666667
$ws.Range("H3").Value = 6.427913
$ws.Range("I3").Value = 0.05953067687027142
$ws.Range("J3").Value = 0.05953067687027141
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 160.8390196666666
$ws.Range("N3").Value = 482.517059
$ws.Range("O3").Value = 0.9948614511421033
$ws.Range("P3").Value = 0.9948614511421032
$ws.Range("Q3").Value = 344.6197418075408
$ws.Range("R3").Value = 3101.577676267867
$ws.Range("S3").Value = 0.05922477557862987
$ws.Range("T3").Value = 0.05922477557862985
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.142637666666667
$ws.Range("H4").Value = 6.427913
$ws.Range("I4").Value = 0.05953067687027142
$ws.Range("J4").Value = 0.05953067687027141
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.119972
$ws.Range("N4").Value = 0.359916
$ws.Range("O4").Value = 0.0007420806111836584
$ws.Range("P4").Value = 0.0007420806111836582
$ws.Range("Q4").Value = 0.2570565261453334
$ws.Range("R4").Value = 2.313508735308
$ws.Range("S4").Value = 0.00004417656107606789
$ws.Range("T4").Value = 0.00004417656107606787
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.142637666666667
$ws.Range("H5").Value = 6.427913
$ws.Range("I5").Value = 0.05953067687027142
$ws.Range("J5").Value = 0.05953067687027141
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.474493
$ws.Range("N5").Value = 1.423479
$ws.Range("O5").Value = 0.002934951950808252
$ws.Range("P5").Value = 0.002934951950808251
$ws.Range("Q5").Value = 1.016666574369667
$ws.Range("R5").Value = 9.149999169327
$ws.Range("S5").Value = 0.0001747196762133388
$ws.Range("T5").Value = 0.0001747196762133387
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.734454666666666
$ws.Range("H6").Value = 14.203364
$ws.Range("I6").Value = 0.1315412751782492
$ws.Range("J6").Value = 0.1315412751782492
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.236283
$ws.Range("N6").Value = 0.708849
$ws.Range("O6").Value = 0.001461516295904947
$ws.Range("P6").Value = 0.001461516295904947
$ws.Range("Q6").Value = 1.118671152004
$ws.Range("R6").Value = 10.068040368036
$ws.Range("S6").Value = 0.0001922497172571282
$ws.Range("T6").Value = 0.0001922497172571281
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.734454666666666
$ws.Range("H7").Value = 14.203364
$ws.Range("I7").Value = 0.1315412751782492
$ws.Range("J7").Value = 0.1315412751782492
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 160.8390196666666
$ws.Range("N7").Value = 482.517059
$ws.Range("O7").Value = 0.9948614511421033
$ws.Range("P7").Value = 0.9948614511421032
$ws.Range("Q7").Value = 761.4850472429416
$ws.Range("R7").Value = 6853.365425186475
$ws.Range("S7").Value = 0.1308653439089158
$ws.Range("T7").Value = 0.1308653439089157
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.734454666666666
$ws.Range("H8").Value = 14.203364
$ws.Range("I8").Value = 0.1315412751782492
$ws.Range("J8").Value = 0.1315412751782492
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.119972
$ws.Range("N8").Value = 0.359916
$ws.Range("O8").Value = 0.0007420806111836584
$ws.Range("P8").Value = 0.0007420806111836582
$ws.Range("Q8").Value = 0.5680019952693334
$ws.Range("R8").Value = 5.112017957423999
$ws.Range("S8").Value = 0.00009761422988015298
$ws.Range("T8").Value = 0.00009761422988015294
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.734454666666666
$ws.Range("H9").Value = 14.203364
$ws.Range("I9").Value = 0.1315412751782492
$ws.Range("J9").Value = 0.1315412751782492
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.474493
$ws.Range("N9").Value = 1.423479
$ws.Range("O9").Value = 0.002934951950808252
$ws.Range("P9").Value = 0.002934951950808251
$ws.Range("Q9").Value = 2.246465598150666
$ws.Range("R9").Value = 20.218190383356
$ws.Range("S9").Value = 0.0003860673221962077
$ws.Range("T9").Value = 0.0003860673221962075
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.385605666666667
$ws.Range("H10").Value = 13.156817
$ws.Range("I10").Value = 0.1218489144872206
$ws.Range("J10").Value = 0.1218489144872206
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.236283
$ws.Range("N10").Value = 0.708849
$ws.Range("O10").Value = 0.001461516295904947
$ws.Range("P10").Value = 0.001461516295904947
$ws.Range("Q10").Value = 1.036244063737
$ws.Range("R10").Value = 9.326196573633
$ws.Range("S10").Value = 0.0001780841741614013
$ws.Range("T10").Value = 0.0001780841741614013
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.385605666666667
$ws.Range("H11").Value = 13.156817
$ws.Range("I11").Value = 0.1218489144872206
$ws.Range("J11").Value = 0.1218489144872206
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 160.8390196666666
$ws.Range("N11").Value = 482.517059
$ws.Range("O11").Value = 0.9948614511421033
$ws.Range("P11").Value = 0.9948614511421032
$ws.Range("Q11").Value = 705.3765160712447
$ws.Range("R11").Value = 6348.388644641203
$ws.Range("S11").Value = 0.1212227878868464
$ws.Range("T11").Value = 0.1212227878868463
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.385605666666667
$ws.Range("H12").Value = 13.156817
$ws.Range("I12").Value = 0.1218489144872206
$ws.Range("J12").Value = 0.1218489144872206
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.119972
$ws.Range("N12").Value = 0.359916
$ws.Range("O12").Value = 0.0007420806111836584
$ws.Range("P12").Value = 0.0007420806111836582
$ws.Range("Q12").Value = 0.5261498830413334
$ws.Range("R12").Value = 4.735348947372001
$ws.Range("S12").Value = 0.000090421716934742
$ws.Range("T12").Value = 0.00009042171693474195
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.385605666666667
$ws.Range("H13").Value = 13.156817
$ws.Range("I13").Value = 0.1218489144872206
$ws.Range("J13").Value = 0.1218489144872206
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.474493
$ws.Range("N13").Value = 1.423479
$ws.Range("O13").Value = 0.002934951950808252
$ws.Range("P13").Value = 0.002934951950808251
$ws.Range("Q13").Value = 2.080939189593667
$ws.Range("R13").Value = 18.728452706343
$ws.Range("S13").Value = 0.000357620709278136
$ws.Range("T13").Value = 0.0003576207092781359
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 24.72946233333333
$ws.Range("H14").Value = 74.188387
$ws.Range("I14").Value = 0.6870791334642589
$ws.Range("J14").Value = 0.6870791334642588
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.236283
$ws.Range("N14").Value = 0.708849
$ws.Range("O14").Value = 0.001461516295904947
$ws.Range("P14").Value = 0.001461516295904947
$ws.Range("Q14").Value = 5.843151548507
$ws.Range("R14").Value = 52.588363936563
$ws.Range("S14").Value = 0.001004177350134264
$ws.Range("T14").Value = 0.001004177350134264
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 24.72946233333333
$ws.Range("H15").Value = 74.188387
$ws.Range("I15").Value = 0.6870791334642589
$ws.Range("J15").Value = 0.6870791334642588
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 160.8390196666666
$ws.Range("N15").Value = 482.517059
$ws.Range("O15").Value = 0.9948614511421033
$ws.Range("P15").Value = 0.9948614511421032
$ws.Range("Q15").Value = 3977.462478577092
$ws.Range("R15").Value = 35797.16230719383
$ws.Range("S15").Value = 0.6835485437677115
$ws.Range("T15").Value = 0.6835485437677112
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 24.72946233333333
$ws.Range("H16").Value = 74.188387
$ws.Range("I16").Value = 0.6870791334642589
$ws.Range("J16").Value = 0.6870791334642588
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.119972
$ws.Range("N16").Value = 0.359916
$ws.Range("O16").Value = 0.0007420806111836584
$ws.Range("P16").Value = 0.0007420806111836582
$ws.Range("Q16").Value = 2.966843055054667
$ws.Range("R16").Value = 26.701587495492
$ws.Range("S16").Value = 0.0005098681032926956
$ws.Range("T16").Value = 0.0005098681032926954
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 24.72946233333333
$ws.Range("H17").Value = 74.188387
$ws.Range("I17").Value = 0.6870791334642589
$ws.Range("J17").Value = 0.6870791334642588
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.474493
$ws.Range("N17").Value = 1.423479
$ws.Range("O17").Value = 0.002934951950808252
$ws.Range("P17").Value = 0.002934951950808251
$ws.Range("Q17").Value = 11.73395677093033
$ws.Range("R17").Value = 105.605610938373
$ws.Range("S17").Value = 0.002016544243120569
$ws.Range("T17").Value = 0.002016544243120569

Write-Output "Applied updated NATMI values to 16 rows (E,G,H,I,J,K,M,N,O,P,Q,R,S,T)"
